# "Fruta / hortaliza, semanal"
#
# The weekly refresh reshuffles which market-day record lands on which
# sheet row. For every row 2..37 the columns Fecha(D), Volumen(J),
# Precio minimo(K), Precio maximo(L), Precio promedio ponderado(M),
# Unidad de comercializacion(N), Origen(O), Precio $/Kg(P) and
# Kg o Unidades(Q) are replaced by the values that (in the prior
# version of the sheet) lived on a different row. Columns A,B,C,E,F,
# G,H,I,R are untouched. Below is the exact row permutation describing
# "new row -> old row it copies its D..Q values from".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices: D=4 J=10 K=11 L=12 M=13 N=14 O=15 P=16 Q=17
$cols = @(4, 10, 11, 12, 13, 14, 15, 16, 17)

# destination row -> source row (values are read from the ORIGINAL sheet
# state, i.e. before any of this script's writes happen)
$rowMap = @(
    @(2, 12), @(3, 16), @(4, 17), @(5, 34), @(6, 25), @(7, 6), @(8, 28),
    @(9, 24), @(10, 7), @(11, 26), @(12, 20), @(13, 19), @(14, 21),
    @(15, 33), @(16, 27), @(17, 13), @(18, 37), @(19, 3), @(20, 8),
    @(21, 18), @(22, 2), @(23, 35), @(24, 29), @(25, 31), @(26, 9),
    @(27, 14), @(28, 23), @(29, 30), @(30, 32), @(31, 36), @(32, 4),
    @(33, 22), @(34, 11), @(35, 5), @(36, 15), @(37, 10)
)

# 1) Snapshot the original D..Q values (by column) for every data row so
#    that later writes never clobber a value still needed as a source.
$snapshot = @{}
for ($r = 2; $r -le 37; $r++) {
    $rowvals = @()
    foreach ($c in $cols) {
        $rowvals += ,($ws.Cells.Item($r, $c).Value())
    }
    $snapshot[$r] = $rowvals
}

# 2) Write the permuted values back out.
foreach ($pair in $rowMap) {
    $destRow = $pair[0]
    $srcRow = $pair[1]
    $vals = $snapshot[$srcRow]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($destRow, $cols[$i]).Value = $vals[$i]
    }
}
